# ---------------------------------------------------------------------------
# Add a "Sort by ID" worksheet that holds a copy of Sheet1's key columns
# (Group, ID, Before dir, After dir, General notes) sorted by ID ascending.
# ---------------------------------------------------------------------------

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# New sheet goes right after Sheet1 and becomes sheet #2.
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$newSheet.Name = "Sort by ID"

# Copy A1:E24 (header + the Group/ID/Before dir/After dir/General notes columns)
# from Sheet1 into the new sheet, preserving cell formatting.
$srcRange = $ws1.Range("A1:E24")
$srcRange.Copy($newSheet.Range("A1"))

# Sort the copied data (rows 2-24) by column B (ID) ascending. Using the
# worksheet's persistent Sort object (rather than a one-shot Range.Sort)
# so the workbook records the sort state, matching a real "Data > Sort".
$sortDataRange = $newSheet.Range("A2:E24")
$sortKeyRange  = $newSheet.Range("B2:B24")

$newSheet.Sort.SortFields.Clear()
$newSheet.Sort.SortFields.Add($sortKeyRange)
$newSheet.Sort.SetRange($sortDataRange)
$newSheet.Sort.Header = 2
$newSheet.Sort.Apply()

# Tidy up column B's width on the new sheet (best-fit for the ID numbers).
$newSheet.Columns.Item(2).ColumnWidth = 10.17

# Restore/update the on-screen selections: Sheet1 keeps A1:E24 highlighted,
# the new sheet opens with column B selected, and the new sheet becomes
# the active tab (as it was right after being created/sorted).
$ws1.Range("A1:E24").Select()
$newSheet.Range("B1:B1048576").Select()
$newSheet.Activate()
